# Update the table style applied to every table in the deck from the
# custom "Table_0" style ({4E98B519-F85A-4075-9106-43047619D924}) to the
# built-in table style {3963827C-B6BB-41C1-A7B3-501B7965716B}.
#
# Three slides in this deck (the tables comparing business structures)
# carry a table each; walk every slide/shape and re-apply the style to
# any shape that actually hosts a table so the edit is robust even if
# shape indices shift.

$p = $ppt.ActivePresentation
$newStyleId = "{3963827C-B6BB-41C1-A7B3-501B7965716B}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
